$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, matching the left
# neighbour's (M) width - mirrors what Excel does when a column is
# inserted via the UI ("Insert" copies formatting from the column to
# the left, shifting the old N/O/P columns one place to the right).
$mWidth = $ws.Columns("M:M").ColumnWidth
$ws.Columns("N:N").Insert() | Out-Null
$ws.Columns("N:N").ColumnWidth = $mWidth

# Switch focus to the "Repayment schedule" sheet and leave the
# selection on K17, matching the recorded UI state (this also moves
# tabSelected/activeTab off the NewLoanInput sheet and onto this one).
$ws.Activate() | Out-Null
$ws.Range("K17").Select() | Out-Null
